# Scheduled data refresh: update market-board price/profit figures on each
# job sheet (currentAveragePrice* / LevePrice* / LeveProfit* columns).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road
$ws.Range("H17").Value = 2778641.2
$ws.Range("J17").Value = 2778641.2
$ws.Range("L17").Value = 8335923.600000001
$ws.Range("N17").Value = -8336259.600000001
# Row 28: The Writing Is Not on the Wall
$ws.Range("H28").Value = 728.55554
$ws.Range("I28").Value = 708.4286
$ws.Range("J28").Value = 799
$ws.Range("K28").Value = 708.4286
$ws.Range("L28").Value = 799
$ws.Range("M28").Value = -223.4286
$ws.Range("N28").Value = -1769
# Row 112: Making Ends Meet
$ws.Range("H112").Value = 1061.2084
$ws.Range("J112").Value = 1090.1818
$ws.Range("L112").Value = 3270.5454
$ws.Range("N112").Value = -5486.5454
# Row 113: Amaro Kart
$ws.Range("H113").Value = 7086.933
$ws.Range("I113").Value = 2780.5334
$ws.Range("J113").Value = 11393.333
$ws.Range("K113").Value = 2780.5334
$ws.Range("L113").Value = 11393.333
$ws.Range("M113").Value = 473.4666000000002
$ws.Range("N113").Value = -17901.333
# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 31793.516
$ws.Range("I132").Value = 37599.48
$ws.Range("J132").Value = 5666.6665
$ws.Range("K132").Value = 112798.44
$ws.Range("L132").Value = 16999.9995
$ws.Range("M132").Value = -110268.44
$ws.Range("N132").Value = -22059.9995
# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 1549.1316
$ws.Range("I137").Value = 1385.4706
$ws.Range("J137").Value = 1681.619
$ws.Range("K137").Value = 4156.4118
$ws.Range("L137").Value = 5044.857
$ws.Range("M137").Value = -1606.4118
$ws.Range("N137").Value = -10144.857
# Row 138: All-night Crafting
$ws.Range("H138").Value = 2668.4856
$ws.Range("I138").Value = 1585.4375
$ws.Range("J138").Value = 3580.5264
$ws.Range("K138").Value = 4756.3125
$ws.Range("L138").Value = 10741.5792
$ws.Range("M138").Value = 383.6875
$ws.Range("N138").Value = -21021.5792

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 4763.7666
$ws.Range("I32").Value = 3739.4375
$ws.Range("J32").Value = 12958.4
$ws.Range("K32").Value = 3739.4375
$ws.Range("L32").Value = 12958.4
$ws.Range("M32").Value = -3452.4375
$ws.Range("N32").Value = -13532.4
# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 1685.2565
$ws.Range("I61").Value = 1343.5555
$ws.Range("J61").Value = 2454.0833
$ws.Range("K61").Value = 1343.5555
$ws.Range("L61").Value = 2454.0833
$ws.Range("M61").Value = -1131.5555
$ws.Range("N61").Value = -2878.0833
# Row 122: Haste for High Durium
$ws.Range("H122").Value = 1305.3948
$ws.Range("I122").Value = 890.9231
$ws.Range("J122").Value = 2203.4167
$ws.Range("K122").Value = 2672.7693
$ws.Range("L122").Value = 6610.250100000001
$ws.Range("M122").Value = -222.7692999999999
$ws.Range("N122").Value = -11510.2501
# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 2137.1562
$ws.Range("I132").Value = 1330.4706
$ws.Range("J132").Value = 3051.4
$ws.Range("K132").Value = 3991.4118
$ws.Range("L132").Value = 9154.200000000001
$ws.Range("M132").Value = -1461.4118
$ws.Range("N132").Value = -14214.2
# Row 136: Metal with Mettle
$ws.Range("H136").Value = 1685.2565
$ws.Range("I136").Value = 1343.5555
$ws.Range("J136").Value = 2454.0833
$ws.Range("K136").Value = 4030.6665
$ws.Range("L136").Value = 7362.249899999999
$ws.Range("M136").Value = -1480.6665
$ws.Range("N136").Value = -12462.2499

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 3694.8572
$ws.Range("I86").Value = 3799.6
$ws.Range("J86").Value = 3433
$ws.Range("K86").Value = 3799.6
$ws.Range("L86").Value = 3433
$ws.Range("M86").Value = -2676.6
$ws.Range("N86").Value = -5679
# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 3694.8572
$ws.Range("I89").Value = 3799.6
$ws.Range("J89").Value = 3433
$ws.Range("K89").Value = 18998
$ws.Range("L89").Value = 17165
$ws.Range("M89").Value = -13382
$ws.Range("N89").Value = -28397
# Row 94: High Steal
$ws.Range("H94").Value = 1593.8
$ws.Range("I94").Value = 1742.25
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 1742.25
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = -1291.25
$ws.Range("N94").Value = -1902
# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 2913.3333
$ws.Range("I105").Value = 2784.2856
$ws.Range("J105").Value = 3094
$ws.Range("K105").Value = 2784.2856
$ws.Range("L105").Value = 3094
$ws.Range("M105").Value = -1037.2856
$ws.Range("N105").Value = -6588
# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 1371.0851
$ws.Range("I134").Value = 943.4167
$ws.Range("J134").Value = 2770.7273
$ws.Range("K134").Value = 2830.2501
$ws.Range("L134").Value = 8312.1819
$ws.Range("M134").Value = -295.2501000000002
$ws.Range("N134").Value = -13382.1819

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof
$ws.Range("H16").Value = 911.4286
$ws.Range("I16").Value = 896.6667
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 896.6667
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -609.6667
$ws.Range("N16").Value = -1574
# Row 113: Patient Patients
$ws.Range("H113").Value = 911.4286
$ws.Range("I113").Value = 896.6667
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 896.6667
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1273.3333
$ws.Range("N113").Value = -5340

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap
$ws.Range("H5").Value = 702598.8
$ws.Range("I5").Value = 414.8095
$ws.Range("K5").Value = 1244.4285
$ws.Range("M5").Value = -1132.4285
# Row 80: Saucy for a Suitor
$ws.Range("H80").Value = 4632.75
$ws.Range("J80").Value = 4979.3
$ws.Range("L80").Value = 14937.9
$ws.Range("N80").Value = -16809.9
# Row 83: Saved by the Sauce (L)
$ws.Range("H83").Value = 4632.75
$ws.Range("J83").Value = 4979.3
$ws.Range("L83").Value = 44813.7
$ws.Range("N83").Value = -54173.7
# Row 132: More Mezcal
$ws.Range("H132").Value = 1575.1666
$ws.Range("J132").Value = 1733.3334
$ws.Range("L132").Value = 15600.0006
$ws.Range("N132").Value = -20660.0006
# Row 135: Not-so-secret Ingredient
$ws.Range("H135").Value = 702598.8
$ws.Range("I135").Value = 414.8095
$ws.Range("K135").Value = 3733.2855
$ws.Range("M135").Value = -1198.2855

$ws = $wb.Worksheets.Item("GSM")
# Row 132: On Board for Lar
$ws.Range("H132").Value = 3029
$ws.Range("I132").Value = 2715.1667
$ws.Range("K132").Value = 8145.500100000001
$ws.Range("M132").Value = -5615.500100000001

$ws = $wb.Worksheets.Item("LTW")
# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 5621.3213
$ws.Range("I132").Value = 7065.0713
$ws.Range("J132").Value = 4177.5713
$ws.Range("K132").Value = 21195.2139
$ws.Range("L132").Value = 12532.7139
$ws.Range("M132").Value = -18665.2139
$ws.Range("N132").Value = -17592.7139
# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 6310111
$ws.Range("I136").Value = 32870.97
$ws.Range("J136").Value = 15875429
$ws.Range("K136").Value = 98612.91
$ws.Range("L136").Value = 47626287
$ws.Range("M136").Value = -96062.91
$ws.Range("N136").Value = -47631387

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins
$ws.Range("H132").Value = 1402.8302
$ws.Range("I132").Value = 958.32556
$ws.Range("K132").Value = 2874.97668
$ws.Range("M132").Value = -344.9766799999998
# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 5645.8
$ws.Range("I136").Value = 905
$ws.Range("J136").Value = 12757
$ws.Range("K136").Value = 2715
$ws.Range("L136").Value = 38271
$ws.Range("M136").Value = -165
$ws.Range("N136").Value = -43371
